# Apply the "LinuxForHealth" re-branding / version bump edits to the
# StructureDefinition-claim-status workbook.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-status"
# Version
$meta.Range("B3").Value = "8.0.0"
# Date
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
# Publisher
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet ---------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Row 2 = "Extension" : the (duplicated) FHIR invariant text is cleared
# from this row -- it remains correctly on row 4 ("Extension.extension").
$elements.Range("AI2").Value = ""

# Row 5 = "Extension.url" : Fixed Value mirrors the updated URL
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/claim-status"
